# Update the two worksheets ("展览" and "全部类型") that contain the
# "南宁·2024良牙动漫秋季盛典（秋典）" (row 5) and
# "南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini" (row 6) entries.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # F5: 想去人数 (want-to-go count) 3210 -> 3228
    $ws.Range("F5").Value = 3228

    # I5: Cover image URL update
    $ws.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202409/ALSTkhdX1725520827309.jpeg"

    # F6: 想去人数 (want-to-go count) 325 -> 327
    $ws.Range("F6").Value = 327
}
